# Applies the September edits to the Vietnam Town Condominium Owners
# Association courtesy-notice letter:
#   1. Bumps the letter date from September 19, 2025 to September 21, 2025.
#   2. Splits the single-line mailing address "999 Story Road, San Jose CA
#      95122" into two paragraphs: "999 Story Road" and "San Jose, CA 95122".
#   3. Removes the now-superfluous blank "No Spacing" paragraph that used to
#      sit directly under "...Board of Directors".

$d = $word.ActiveDocument

# 1. Letter date.
$d.Content.Find.Execute("September 19, 2025", $true, $false, $false, $false, `
    $false, $true, 1, $false, "September 21, 2025", 2) | Out-Null

# 2. Split the mailing-address paragraph into two paragraphs, preserving the
#    paragraph/run formatting (Arial 11pt, autoSpaceDE/DN off) on both halves.
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -eq "999 Story Road, San Jose CA 95122`r") {
        $p.Range.Text = "999 Story Road`rSan Jose, CA 95122"
        break
    }
}

# 3. Drop the empty "No Spacing" paragraph right after "...Board of
#    Directors" (the signature block now runs straight into the Title-style
#    spacer paragraphs that follow it).
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -eq "Vietnam Town Condominium Owners Association Board of Directors`r") {
        $next = $d.Paragraphs.Item($i + 1)
        if ($next.Range.Text -eq "`r" -and $next.Style.NameLocal -eq "No Spacing") {
            $next.Range.Delete()
        }
        break
    }
}
